$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.330421
$ws.Range("N2").Value = 0.991263
$ws.Range("O2").Value = 0.1900853910865743
$ws.Range("P2").Value = 0.1900853910865742
$ws.Range("Q2").Value = 0.4227715768336667
$ws.Range("R2").Value = 3.804944191503
$ws.Range("S2").Value = 0.1900853910865743
$ws.Range("T2").Value = 0.1900853910865742

# Row 3
$ws.Range("O3").Value = 0.515657077987202
$ws.Range("P3").Value = 0.515657077987202
$ws.Range("S3").Value = 0.515657077987202
$ws.Range("T3").Value = 0.515657077987202

# Row 4
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.511501
$ws.Range("N4").Value = 1.534503
$ws.Range("O4").Value = 0.2942575309262239
$ws.Range("P4").Value = 0.2942575309262239
$ws.Range("Q4").Value = 0.6544622899936666
$ws.Range("R4").Value = 5.890160609943
$ws.Range("S4").Value = 0.2942575309262239
$ws.Range("T4").Value = 0.2942575309262239
